$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibition) sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1710
$wsExhibit.Range("F3").Value = 7879

# Update "全部类型" (All Types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1710
$wsAll.Range("F3").Value = 7879
